$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 660.2353000000001
$ws.Cells.Item(41, 9).Value = 76
$ws.Cells.Item(41, 10).Value = 903.6667
$ws.Cells.Item(41, 11).Value = 76
$ws.Cells.Item(41, 12).Value = 903.6667
$ws.Cells.Item(41, 13).Value = 364
$ws.Cells.Item(41, 14).Value = -1783.6667
$ws.Cells.Item(53, 8).Value = 2512.6155
$ws.Cells.Item(53, 9).Value = 184.16667
$ws.Cells.Item(53, 10).Value = 4508.4287
$ws.Cells.Item(53, 11).Value = 184.16667
$ws.Cells.Item(53, 12).Value = 4508.4287
$ws.Cells.Item(53, 13).Value = 452.83333
$ws.Cells.Item(53, 14).Value = -5782.4287
$ws.Cells.Item(97, 8).Value = 986.913
$ws.Cells.Item(97, 10).Value = 1022.6818
$ws.Cells.Item(97, 12).Value = 3068.0454
$ws.Cells.Item(97, 14).Value = -4060.0454
$ws.Cells.Item(113, 8).Value = 50004136
$ws.Cells.Item(113, 10).Value = 6899.4
$ws.Cells.Item(113, 12).Value = 6899.4
$ws.Cells.Item(113, 14).Value = -13407.4
$ws.Cells.Item(116, 8).Value = 4266.9443
$ws.Cells.Item(116, 9).Value = 2240.9
$ws.Cells.Item(116, 11).Value = 2240.9
$ws.Cells.Item(116, 13).Value = 1201.1
$ws.Cells.Item(129, 8).Value = 164845.12
$ws.Cells.Item(129, 10).Value = 189693.8
$ws.Cells.Item(129, 12).Value = 569081.3999999999
$ws.Cells.Item(129, 14).Value = -579081.3999999999
$ws.Cells.Item(132, 8).Value = 3970.64
$ws.Cells.Item(132, 9).Value = 4073.3333
$ws.Cells.Item(132, 10).Value = 1506
$ws.Cells.Item(132, 11).Value = 12219.9999
$ws.Cells.Item(132, 12).Value = 4518
$ws.Cells.Item(132, 13).Value = -9689.999899999999
$ws.Cells.Item(132, 14).Value = -9578
$ws.Cells.Item(137, 8).Value = 1187.5918
$ws.Cells.Item(137, 9).Value = 1199.9459
$ws.Cells.Item(137, 11).Value = 3599.8377
$ws.Cells.Item(137, 13).Value = -1049.8377
$ws.Cells.Item(138, 8).Value = 2430.75
$ws.Cells.Item(138, 10).Value = 3457.9443
$ws.Cells.Item(138, 12).Value = 10373.8329
$ws.Cells.Item(138, 14).Value = -20653.8329

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 558.2368
$ws.Cells.Item(2, 9).Value = 589.65515
$ws.Cells.Item(2, 10).Value = 457
$ws.Cells.Item(2, 11).Value = 589.65515
$ws.Cells.Item(2, 12).Value = 457
$ws.Cells.Item(2, 13).Value = -476.65515
$ws.Cells.Item(2, 14).Value = -683
$ws.Cells.Item(32, 8).Value = 6478.4614
$ws.Cells.Item(32, 9).Value = 4473.932
$ws.Cells.Item(32, 10).Value = 12703.053
$ws.Cells.Item(32, 11).Value = 4473.932
$ws.Cells.Item(32, 12).Value = 12703.053
$ws.Cells.Item(32, 13).Value = -4186.932
$ws.Cells.Item(32, 14).Value = -13277.053
$ws.Cells.Item(45, 8).Value = 2858.1
$ws.Cells.Item(45, 9).Value = 2019.8
$ws.Cells.Item(45, 11).Value = 2019.8
$ws.Cells.Item(45, 13).Value = -1642.8
$ws.Cells.Item(61, 8).Value = 2236.348
$ws.Cells.Item(61, 9).Value = 1613.8823
$ws.Cells.Item(61, 10).Value = 4000
$ws.Cells.Item(61, 11).Value = 1613.8823
$ws.Cells.Item(61, 12).Value = 4000
$ws.Cells.Item(61, 13).Value = -1401.8823
$ws.Cells.Item(61, 14).Value = -4424
$ws.Cells.Item(74, 8).Value = 33334906
$ws.Cells.Item(74, 9).Value = 47619504
$ws.Cells.Item(74, 10).Value = 4179.222
$ws.Cells.Item(74, 11).Value = 47619504
$ws.Cells.Item(74, 12).Value = 4179.222
$ws.Cells.Item(74, 13).Value = -47618630
$ws.Cells.Item(74, 14).Value = -5927.222
$ws.Cells.Item(77, 8).Value = 33334906
$ws.Cells.Item(77, 9).Value = 47619504
$ws.Cells.Item(77, 10).Value = 4179.222
$ws.Cells.Item(77, 11).Value = 238097520
$ws.Cells.Item(77, 12).Value = 20896.11
$ws.Cells.Item(77, 13).Value = -238093152
$ws.Cells.Item(77, 14).Value = -29632.11
$ws.Cells.Item(116, 8).Value = 558.2368
$ws.Cells.Item(116, 9).Value = 589.65515
$ws.Cells.Item(116, 10).Value = 457
$ws.Cells.Item(116, 11).Value = 589.65515
$ws.Cells.Item(116, 12).Value = 457
$ws.Cells.Item(116, 13).Value = 1704.34485
$ws.Cells.Item(116, 14).Value = -5045
$ws.Cells.Item(136, 8).Value = 2236.348
$ws.Cells.Item(136, 9).Value = 1613.8823
$ws.Cells.Item(136, 10).Value = 4000
$ws.Cells.Item(136, 11).Value = 4841.6469
$ws.Cells.Item(136, 12).Value = 12000
$ws.Cells.Item(136, 13).Value = -2291.6469
$ws.Cells.Item(136, 14).Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 558.2368
$ws.Cells.Item(3, 9).Value = 589.65515
$ws.Cells.Item(3, 10).Value = 457
$ws.Cells.Item(3, 11).Value = 589.65515
$ws.Cells.Item(3, 12).Value = 457
$ws.Cells.Item(3, 13).Value = -475.65515
$ws.Cells.Item(3, 14).Value = -685
$ws.Cells.Item(25, 8).Value = 484
$ws.Cells.Item(25, 9).Value = 484
$ws.Cells.Item(25, 11).Value = 484
$ws.Cells.Item(25, 13).Value = -249
$ws.Cells.Item(107, 8).Value = 821.44446
$ws.Cells.Item(107, 9).Value = 679.06665
$ws.Cells.Item(107, 10).Value = 1533.3334
$ws.Cells.Item(107, 11).Value = 679.06665
$ws.Cells.Item(107, 12).Value = 1533.3334
$ws.Cells.Item(107, 13).Value = 1240.93335
$ws.Cells.Item(107, 14).Value = -5373.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 468.7143
$ws.Cells.Item(22, 10).Value = 594.2
$ws.Cells.Item(22, 12).Value = 594.2
$ws.Cells.Item(22, 14).Value = -1294.2
$ws.Cells.Item(31, 8).Value = 4129.1377
$ws.Cells.Item(31, 10).Value = 4242.579
$ws.Cells.Item(31, 12).Value = 4242.579
$ws.Cells.Item(31, 14).Value = -4832.579
$ws.Cells.Item(34, 8).Value = 4129.1377
$ws.Cells.Item(34, 10).Value = 4242.579
$ws.Cells.Item(34, 12).Value = 4242.579
$ws.Cells.Item(34, 14).Value = -4646.579
$ws.Cells.Item(58, 8).Value = 22475.791
$ws.Cells.Item(58, 9).Value = 1659.8462
$ws.Cells.Item(58, 11).Value = 1659.8462
$ws.Cells.Item(58, 13).Value = -1456.8462
$ws.Cells.Item(62, 8).Value = 55559456
$ws.Cells.Item(62, 9).Value = 71431944
$ws.Cells.Item(62, 11).Value = 71431944
$ws.Cells.Item(62, 13).Value = -71431320
$ws.Cells.Item(65, 8).Value = 55559456
$ws.Cells.Item(65, 9).Value = 71431944
$ws.Cells.Item(65, 11).Value = 357159720
$ws.Cells.Item(65, 13).Value = -357156600
$ws.Cells.Item(132, 8).Value = 5051
$ws.Cells.Item(132, 9).Value = 4117.5713
$ws.Cells.Item(132, 10).Value = 6140
$ws.Cells.Item(132, 11).Value = 12352.7139
$ws.Cells.Item(132, 12).Value = 18420
$ws.Cells.Item(132, 13).Value = -9822.713899999999
$ws.Cells.Item(132, 14).Value = -23480
$ws.Cells.Item(134, 8).Value = 1775
$ws.Cells.Item(134, 9).Value = 1825
$ws.Cells.Item(134, 11).Value = 5475
$ws.Cells.Item(134, 13).Value = -2940
$ws.Cells.Item(136, 8).Value = 22475.791
$ws.Cells.Item(136, 9).Value = 1659.8462
$ws.Cells.Item(136, 11).Value = 4979.5386
$ws.Cells.Item(136, 13).Value = -2429.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 855
$ws.Cells.Item(23, 9).Value = 20
$ws.Cells.Item(23, 10).Value = 1133.3334
$ws.Cells.Item(23, 11).Value = 60
$ws.Cells.Item(23, 12).Value = 3400.0002
$ws.Cells.Item(23, 13).Value = 175
$ws.Cells.Item(23, 14).Value = -3870.0002
$ws.Cells.Item(131, 8).Value = 707.34
$ws.Cells.Item(131, 10).Value = 707.34
$ws.Cells.Item(131, 12).Value = 2122.02
$ws.Cells.Item(131, 14).Value = -12202.02

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 20001642
$ws.Cells.Item(102, 9).Value = 22728544
$ws.Cells.Item(102, 10).Value = 4354.6665
$ws.Cells.Item(102, 11).Value = 22728544
$ws.Cells.Item(102, 12).Value = 4354.6665
$ws.Cells.Item(102, 13).Value = -22726922
$ws.Cells.Item(102, 14).Value = -7598.6665
$ws.Cells.Item(107, 8).Value = 3496778.5
$ws.Cells.Item(107, 9).Value = 293.46155
$ws.Cells.Item(107, 10).Value = 8547257
$ws.Cells.Item(107, 11).Value = 293.46155
$ws.Cells.Item(107, 12).Value = 8547257
$ws.Cells.Item(107, 13).Value = 1626.53845
$ws.Cells.Item(107, 14).Value = -8551097

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2774.0454
$ws.Cells.Item(22, 9).Value = 3805.4
$ws.Cells.Item(22, 11).Value = 3805.4
$ws.Cells.Item(22, 13).Value = -3510.4
$ws.Cells.Item(27, 8).Value = 2774.0454
$ws.Cells.Item(27, 9).Value = 3805.4
$ws.Cells.Item(27, 11).Value = 3805.4
$ws.Cells.Item(27, 13).Value = -3698.4
$ws.Cells.Item(46, 8).Value = 932.4912
$ws.Cells.Item(46, 9).Value = 931.7692
$ws.Cells.Item(46, 10).Value = 940
$ws.Cells.Item(46, 11).Value = 931.7692
$ws.Cells.Item(46, 12).Value = 940
$ws.Cells.Item(46, 13).Value = -743.7692
$ws.Cells.Item(46, 14).Value = -1316
$ws.Cells.Item(111, 8).Value = 29987
$ws.Cells.Item(111, 10).Value = 29987
$ws.Cells.Item(111, 12).Value = 29987
$ws.Cells.Item(111, 14).Value = -38167

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1518.409
$ws.Cells.Item(126, 9).Value = 1387.875
$ws.Cells.Item(126, 10).Value = 1866.5
$ws.Cells.Item(126, 11).Value = 4163.625
$ws.Cells.Item(126, 12).Value = 5599.5
$ws.Cells.Item(126, 13).Value = -1693.625
$ws.Cells.Item(126, 14).Value = -10539.5
$ws.Cells.Item(132, 8).Value = 1387.2963
$ws.Cells.Item(132, 9).Value = 884.5
$ws.Cells.Item(132, 10).Value = 3599.6
$ws.Cells.Item(132, 11).Value = 2653.5
$ws.Cells.Item(132, 12).Value = 10798.8
$ws.Cells.Item(132, 13).Value = -123.5
$ws.Cells.Item(132, 14).Value = -15858.8
$ws.Cells.Item(136, 8).Value = 25179764
$ws.Cells.Item(136, 9).Value = 31281758
$ws.Cells.Item(136, 11).Value = 93845274
$ws.Cells.Item(136, 13).Value = -93842724
